$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J, using the same style as existing headers (style index 1 -> bold/centered/bordered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-10
$values = @(
    @(3, 3),
    @(7, 8),
    @(10, 10),
    @(6, 7),
    @(9, 9),
    @(4, 4),
    @(9, 9),
    @(8, 8),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
